$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Semestre ideal" row (B9/C9 share the same text)
$ws.Range("B9").Value = "EA-4,EB-5,EQD-4,EQN-5"
$ws.Range("C9").Value = "EA-4,EB-5,EQD-4,EQN-5"

# Update existing requirement row 24 text
$ws.Range("B24").Value = "LOB1024 -  Mecânica  (Requisito fraco)`n"
$ws.Range("C24").Value = "LOB1024 -  Mecânica  (Requisito fraco)`n"

# Add new requirement row 25
$ws.Range("B25").Value = "LOB1052 -  Cálculo III  (Requisito fraco)`n"
$ws.Range("C25").Value = "LOB1052 -  Cálculo III  (Requisito fraco)`n"

# Copy formatting from row 24 to row 25 so new row matches existing style
$ws.Range("B24:C24").Copy()
$ws.Range("B25:C25").PasteSpecial(-4122)  # xlPasteFormats

$ws.Rows("25").RowHeight = 30
